$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G5 value
$ws.Range("G5").Value = 54.6

# Update G8:G18 values to 51
$ws.Range("G8:G18").Value = 51

# Update the active cell / selection to G5
$ws.Range("G5").Select()
